$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Column R: fill in pricing_interest_rate_type values (previously blank placeholders) ---
$ws.Range("R1").Value = "pricing_interest_rate_type"
$ws.Range("R1").WrapText = $true

$ws.Range("R2").Value = 1
$ws.Range("R2").WrapText = $true

$ws.Range("R3").Value = 2
$ws.Range("R3").WrapText = $true

$ws.Range("R4").Value = 2
$ws.Range("R4").WrapText = $true

$ws.Range("R5").Value = 3
$ws.Range("R5").WrapText = $true

$ws.Range("R6").Value = 4
$ws.Range("R6").WrapText = $true

$ws.Range("R7").Value = 5
$ws.Range("R7").WrapText = $true

$ws.Range("R8").Value = 6
$ws.Range("R8").WrapText = $true

$ws.Range("R9").Value = 7
$ws.Range("R9").WrapText = $true

$ws.Range("R10").Value = 8
$ws.Range("R10").WrapText = $true

$ws.Range("R11").Value = 9
$ws.Range("R11").WrapText = $true

# --- Column S: new "pricing_fixed_rate" column ---
$ws.Range("S1").Value = "pricing_fixed_rate"
$ws.Range("S1").WrapText = $true

$ws.Range("S2").Style = "Normal"
$ws.Range("S2").Value = "abc123"

$ws.Range("S3").Style = "Normal"
$ws.Range("S3").Value = "abc123"

$ws.Range("S4").Style = "Normal"
$ws.Range("S4").Value = 0

$ws.Range("S5").Style = "Normal"
$ws.Range("S5").Value = 0.1

$ws.Range("S6").Style = "Normal"
$ws.Range("S6").Value = 0.2

# Row 7 in column S is intentionally left untouched (no cell entry at all)

$ws.Range("S8").Style = "Normal"
$ws.Range("S9").Style = "Normal"
$ws.Range("S10").Style = "Normal"
$ws.Range("S11").Style = "Normal"

# --- Update selection to reflect the newly added column ---
$ws.Range("S1:S11").Select()
